# Update "想去人数" (F) counts and cancel three cv meet-and-greet events
# (rows 20-22) on both the "展览" and "全部类型" sheets, which carry the
# same event data.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Simple numeric bumps in column F ("想去人数")
    $ws.Range("F2").Value  = 851
    $ws.Range("F4").Value  = 2163
    $ws.Range("F6").Value  = 12580
    $ws.Range("F10").Value = 458
    $ws.Range("F11").Value = 1145
    $ws.Range("F12").Value = 944
    $ws.Range("F13").Value = 13653
    $ws.Range("F14").Value = 13969
    $ws.Range("F19").Value = 6

    # Row 20: cv 刘圣博 meet-and-greet cancelled
    $ws.Range("C20").Value = "苏州·动漫游戏嘉年华cv刘圣博见面会（取消）"
    $ws.Range("G20").Value = "不可售"

    # Row 21: cv 张文钰 meet-and-greet cancelled
    $ws.Range("C21").Value = "苏州·动漫游戏嘉年华cv张文钰见面会（取消）"
    $ws.Range("F21").Value = 1
    $ws.Range("G21").Value = "不可售"

    # Row 22: cv 沐霏 meet-and-greet cancelled
    $ws.Range("C22").Value = "苏州·动漫游戏嘉年华cv沐霏见面会（取消）"
    $ws.Range("G22").Value = "不可售"

    $ws.Range("F23").Value = 1051
    $ws.Range("F25").Value = 55
    $ws.Range("F26").Value = 537
    $ws.Range("F27").Value = 5080
    $ws.Range("F28").Value = 2
    $ws.Range("F29").Value = 252
}
